# Update "想去人数" (attendance interest count) column F on both the
# "展览" sheet and the "全部类型" sheet, which mirror the same events.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3245
$ws1.Range("F4").Value = 2003
$ws1.Range("F8").Value = 617
$ws1.Range("F15").Value = 10179
$ws1.Range("F20").Value = 8074
$ws1.Range("F33").Value = 7972
$ws1.Range("F37").Value = 85
$ws1.Range("F39").Value = 1435

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3245
$ws4.Range("F6").Value = 2003
$ws4.Range("F12").Value = 617
$ws4.Range("F18").Value = 10179
$ws4.Range("F22").Value = 8074
$ws4.Range("F37").Value = 7972
$ws4.Range("F40").Value = 85
